# Refresh the crypto price/volume snapshot to match the latest GitHub Actions scrape.
# Source cells are plain text (inline strings), including numeric-looking prices like
# "0.999" -- set NumberFormat to Text ("@") before writing those so COM keeps them as
# strings instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textPriceCells = @(
    "D4", "D5", "D6", "D9", "D12", "D14", "D19", "D20", "D21", "D22",
    "D23", "D25", "D29", "D30", "D31", "D32", "D35", "D38", "D39", "D40",
    "D41", "D43", "D44", "D45", "D47", "D48", "D50", "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Bitcoin (row 2)
$ws.Range("D2").Value = '64.756.44'
$ws.Range("E2").Value = '  -0.64%  '

# Ethereum (row 3)
$ws.Range("D3").Value = '3.507.64'
$ws.Range("E3").Value = '  -1.20%  '

# TetherUSD (row 4)
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

# BNB (row 5)
$ws.Range("D5").Value = '587.45'
$ws.Range("E5").Value = '  -1.73%  '

# Solana (row 6)
$ws.Range("D6").Value = '132.58'
$ws.Range("E6").Value = '  -0.69%  '

# LidoStakedEther (row 7)
$ws.Range("D7").Value = '3.506.15'
$ws.Range("E7").Value = '  -1.20%  '

# USDC (row 8)
$ws.Range("E8").Value = '  +0.02%  '

# XRP (row 9)
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -0.90%  '

# Dogecoin (row 10)
$ws.Range("E10").Value = '  +0.91%  '

# Toncoin (row 11)
$ws.Range("E11").Value = '  +0.38%  '

# Cardano (row 12)
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  -0.30%  '

# WrappedliquidstakedEther2.0 (row 13)
$ws.Range("D13").Value = '4.098.47'
$ws.Range("E13").Value = '  -1.38%  '

# Avalanche (row 14)
$ws.Range("D14").Value = '27.78'
$ws.Range("E14").Value = '  +3.24%  '

# ShibaInu (row 15)
$ws.Range("E15").Value = '  -1.39%  '

# TRON (row 16)
$ws.Range("E16").Value = '  +0.66%  '

# WrappedEther (row 17)
$ws.Range("D17").Value = '3.502.70'
$ws.Range("E17").Value = '  -1.34%  '

# WrappedBTC (row 18)
$ws.Range("D18").Value = '64.765.01'
$ws.Range("E18").Value = '  -0.73%  '

# Uniswap (row 19)
$ws.Range("D19").Value = '10.01'
$ws.Range("E19").Value = '  +0.73%  '

# Chainlink (row 20)
$ws.Range("D20").Value = '14.28'
$ws.Range("E20").Value = '  -0.50%  '

# Polkadot (row 21)
$ws.Range("D21").Value = '5.69'
$ws.Range("E21").Value = '  -2.22%  '

# BitcoinCash (row 22)
$ws.Range("D22").Value = '391.41'
$ws.Range("E22").Value = '  +0.28%  '

# Polygon (row 23)
$ws.Range("D23").Value = '0.577'
$ws.Range("E23").Value = '  +0.09%  '

# WrappedeETH (row 24)
$ws.Range("D24").Value = '3.646.64'
$ws.Range("E24").Value = '  -1.24%  '

# Litecoin (row 25)
$ws.Range("D25").Value = '74.17'
$ws.Range("E25").Value = '  +0.23%  '

# Dai (row 26)
$ws.Range("E26").Value = '  +0.11%  '

# PEPE (row 27)
$ws.Range("E27").Value = '  -4.01%  '

# Fetch.AI (row 28)
$ws.Range("E28").Value = '  +1.66%  '

# RenderToken (row 29)
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.46'
$ws.Range("E29").Value = '  -4.41%  '

# Binance-PegBSC-USD (row 30)
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.07%  '

# PancakeSwap (row 31)
$ws.Range("D31").Value = '2.27'
$ws.Range("E31").Value = '  -0.38%  '

# InternetComputer(DFINITY) (row 32)
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -3.88%  '

# RenzoRestakedETH (row 33)
$ws.Range("D33").Value = '3.508.34'
$ws.Range("E33").Value = '  -1.21%  '

# USDe (row 34)
$ws.Range("E34").Value = '  +0.01%  '

# EthereumClassic (row 35)
$ws.Range("D35").Value = '23.98'
$ws.Range("E35").Value = '  -0.22%  '

# Kaspa (row 36)
$ws.Range("E36").Value = '  -0.19%  '

# ImmutableX (row 37)
$ws.Range("E37").Value = '  +4.02%  '

# NEARProtocol (row 38)
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '5.23'
$ws.Range("E38").Value = '  +4.21%  '

# Monero (row 39)
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '171.96'
$ws.Range("E39").Value = '  +0.90%  '

# Aptos (row 40)
$ws.Range("D40").Value = '6.97'
$ws.Range("E40").Value = '  +0.64%  '

# Hedera (row 41)
$ws.Range("D41").Value = '0.0812'
$ws.Range("E41").Value = '  +0.40%  '

# Mantle (row 42)
$ws.Range("E42").Value = '  -1.37%  '

# EnergySwap (row 43)
$ws.Range("D43").Value = '26.22'
$ws.Range("E43").Value = '  -1.09%  '

# FirstDigitalUSD (row 44)
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.05%  '

# OKB (row 45)
$ws.Range("D45").Value = '42.36'
$ws.Range("E45").Value = '  -1.62%  '

# ONDO (row 46)
$ws.Range("E46").Value = '  -2.20%  '

# Filecoin (row 47)
$ws.Range("D47").Value = '4.41'
$ws.Range("E47").Value = '  -0.66%  '

# Stacks (row 48)
$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  -0.11%  '

# Maker (row 49)
$ws.Range("D49").Value = '2.483.45'
$ws.Range("E49").Value = '  +1.14%  '

# Cosmos (row 50)
$ws.Range("D50").Value = '6.88'
$ws.Range("E50").Value = '  -0.42%  '

# SuiNetwork (row 51)
$ws.Range("D51").Value = '0.907'
$ws.Range("E51").Value = '  +4.05%  '
